$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Resistor block: split R1 (91K) out of the old "56K" group (R7,R16,R1) ---
# Insert a brand-new row 3 for the new 91K / R1 part.
$ws.Rows.Item(3).Insert()
$ws.Range("A3").Value = "91K"
$ws.Range("B3").Value = "R1"
$ws.Range("C3").Value = "R_0805_2012Metric"
$ws.Range("D3").Value = "C137484"

# The former "56K" row (R7,R16,R1) is now row 4 - drop R1 from its designator list.
$ws.Range("B4").Value = "R7,R16"

# --- Capacitor block: two new rows for C11,C12 (330pF) and C13 (4.7pF) ---
# These go right after the last resistor row (330K / R2), which is now row 15.
$ws.Rows.Item(16).Insert()
$ws.Range("A16").Value = "330pF"
$ws.Range("B16").Value = "C11,C12"
$ws.Range("C16").Value = "C_0805_2012Metric"
$ws.Range("D16").Value = "C51207"

$ws.Rows.Item(17).Insert()
$ws.Range("A17").Value = "4.7pF"
$ws.Range("B17").Value = "C13"
$ws.Range("C17").Value = "C_0805_2012Metric"
$ws.Range("D17").Value = "C1820"
